$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 54
$ws.Range("B10").Value = "reset_count"
$ws.Range("C10").Value = 1

$ws.Range("G9").Select()
